$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Verifications")

$ws.Range("A12").Value = "LPA Questionnaire"
$ws.Range("A12").Select()
